$wb = $excel.ActiveWorkbook

# --- Sheet "Yearly": correct the June 401K dividend amount ---
$yearly = $wb.Worksheets.Item("Yearly")
$yearly.Range("M8").Value = 58.75
# (downstream totals O8, M15, O15 on this sheet recalc automatically,
#  as do G8/I8/G46/I46 on "All Time" which pull from Yearly!M15 etc.)

# Leave the Yearly sheet's selection where the workbook was last saved
[void]$yearly.Range("M8").Select()

# --- Sheet "All Time": becomes the active/visible sheet, scrolled down ---
$allTime = $wb.Worksheets.Item("All Time")
[void]$allTime.Activate()
$excel.ActiveWindow.ScrollRow = 40
[void]$allTime.Range("M53").Select()
